# Updates cryptos list: prices and 1h volume % changes, plus two row
# re-ranks (ShibaInu/WrappedBTC swap position 15<->16, THORChain -> Mantle at 51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell([string]$addr, [string]$val) {
    $rng = $ws.Range($addr)
    # Force text so numeric-looking strings (prices like "575.66")
    # are not auto-coerced to numbers, then restore the default "Normal"
    # style so no stray number-format style is left on the cell.
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell "D2" '69.689.93'
Set-TextCell "E2" '  +0.25%  '
Set-TextCell "D3" '2.510.13'
Set-TextCell "E3" '  +0.22%  '
Set-TextCell "E4" '  -0.04%  '
Set-TextCell "D5" '575.66'
Set-TextCell "E5" '  +0.25%  '
Set-TextCell "D6" '166.71'
Set-TextCell "E6" '  +0.31%  '
Set-TextCell "E7" '  -0.06%  '
Set-TextCell "E8" '  -0.10%  '
Set-TextCell "D9" '2.509.15'
Set-TextCell "E9" '  +0.16%  '
Set-TextCell "E10" '  +1.97%  '
Set-TextCell "E11" '  -0.36%  '
Set-TextCell "D12" '0.358'
Set-TextCell "E12" '  +5.04%  '
Set-TextCell "D13" '4.94'
Set-TextCell "E13" '  +2.06%  '
Set-TextCell "D14" '2.969.18'
Set-TextCell "E14" '  +0.19%  '
Set-TextCell "B15" 'WrappedBTC'
Set-TextCell "C15" 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextCell "D15" '69.577.54'
Set-TextCell "E15" '  +0.14%  '
Set-TextCell "B16" 'ShibaInu'
Set-TextCell "C16" 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextCell "D16" '0.0000178'
Set-TextCell "E16" '  +2.01%  '
Set-TextCell "D17" '24.91'
Set-TextCell "E17" '  +0.66%  '
Set-TextCell "D18" '2.517.70'
Set-TextCell "E18" '  -0.10%  '
Set-TextCell "D19" '11.29'
Set-TextCell "E19" '  -0.60%  '
Set-TextCell "D20" '7.52'
Set-TextCell "E20" '  -2.72%  '
Set-TextCell "D21" '350.23'
Set-TextCell "E21" '  +0.56%  '
Set-TextCell "D22" '3.94'
Set-TextCell "E22" '  +0.40%  '
Set-TextCell "E23" '  +0.02%  '
Set-TextCell "E24" '  +0.02%  '
Set-TextCell "E25" '  +2.45%  '
Set-TextCell "E26" '  -0.83%  '
Set-TextCell "D27" '8.84'
Set-TextCell "E27" '  -1.25%  '
Set-TextCell "E28" '  -0.05%  '
Set-TextCell "D29" '0.999'
Set-TextCell "E29" '  -0.18%  '
Set-TextCell "D30" '0.0₃0893'
Set-TextCell "E30" '  -0.33%  '
Set-TextCell "D31" '7.84'
Set-TextCell "E31" '  -0.18%  '
Set-TextCell "D32" '461.28'
Set-TextCell "E32" '  -2.41%  '
Set-TextCell "E33" '  -3.71%  '
Set-TextCell "D34" '1.74'
Set-TextCell "E34" '  -0.26%  '
Set-TextCell "E35" '  -0.06%  '
Set-TextCell "D36" '158.94'
Set-TextCell "E36" '  +3.33%  '
Set-TextCell "D37" '0.116'
Set-TextCell "E37" '  +0.72%  '
Set-TextCell "D38" '19.07'
Set-TextCell "E38" '  +0.67%  '
Set-TextCell "D39" '18.48'
Set-TextCell "E39" '  +0.09%  '
Set-TextCell "E40" '  +0.03%  '
Set-TextCell "D41" '0.320'
Set-TextCell "E41" '  +0.74%  '
Set-TextCell "E42" '  -0.53%  '
Set-TextCell "E43" '  +0.16%  '
Set-TextCell "D44" '38.19'
Set-TextCell "E44" '  +0.13%  '
Set-TextCell "D45" '2.22'
Set-TextCell "E45" '  -3.71%  '
Set-TextCell "E46" '  -6.96%  '
Set-TextCell "D47" '142.30'
Set-TextCell "E47" '  -0.64%  '
Set-TextCell "D48" '3.47'
Set-TextCell "E48" '  -1.26%  '
Set-TextCell "E49" '  -1.12%  '
Set-TextCell "D50" '0.0733'
Set-TextCell "E50" '  +0.33%  '
Set-TextCell "B51" 'Mantle'
Set-TextCell "C51" 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextCell "D51" '0.579'
Set-TextCell "E51" '  -1.03%  '
